$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert a new "2022-Q3" row at the top of the "总计" summary sheet,
#    shifting the existing quarters down by one row.
# ---------------------------------------------------------------------
$ws1.Rows.Item(2).Insert()
$ws1.Range("B2:D2").ClearFormats()

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 14
$ws1.Range("D2").Value = 2.09

# Column A keeps the bold / bordered / centered look used by the rest of
# the index column.
$ws1.Range("A2").Font.Bold = $true
$ws1.Range("A2").Borders.LineStyle = 1
$ws1.Range("A2").HorizontalAlignment = -4108
$ws1.Range("A2").VerticalAlignment = -4160

# Re-number the (0-based) index column for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $ws1.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. Add a brand-new "2022-Q3" worksheet (fund holders detail) right
#    after "总计", pushing the other quarter sheets one slot to the
#    right.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q3"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").Borders.LineStyle = 1
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160

$data = @(
    @(0,  "320005", "诺安价值增长混合",               "10.37", "83.41", "8.11", "0.8410", 3),
    @(1,  "003298", "嘉实物流产业股票A",               "9.19",  "83.51", "4.72", "0.4338", 6),
    @(2,  "003299", "嘉实物流产业股票C",               "6.58",  "83.51", "4.72", "0.3106", 6),
    @(3,  "002291", "诺安安鑫灵活配置混合",             "2.66",  "77.38", "6.98", "0.1857", 1),
    @(4,  "013200", "南方均衡优选一年持有期混合A",       "7.27",  "40.60", "1.24", "0.0901", 8),
    @(5,  "001692", "南方国策动力股票",                 "2.60",  "94.22", "2.26", "0.0588", 8),
    @(6,  "012879", "中信建投量化精选6个月持有期混合C", "3.33",  "90.73", "1.22", "0.0406", 4),
    @(7,  "516530", "银华中证现代物流ETF",              "0.89",  "97.53", "3.80", "0.0338", 9),
    @(8,  "516910", "富国中证现代物流ETF",              "0.78",  "99.30", "3.88", "0.0303", 9),
    @(9,  "012426", "南方价值臻选混合A",                "3.91",  "63.50", "0.62", "0.0242", 8),
    @(10, "012878", "中信建投量化精选6个月持有期混合A", "1.67",  "90.73", "1.22", "0.0204", 4),
    @(11, "202213", "南方核心竞争混合",                 "2.05",  "62.82", "0.71", "0.0146", 8),
    @(12, "013201", "南方均衡优选一年持有期混合C",       "0.73",  "40.60", "1.24", "0.0091", 8),
    @(13, "012427", "南方价值臻选混合C",                "0.19",  "63.50", "0.62", "0.0012", 8)
)

$row = 2
foreach ($d in $data) {
    $aCell = $newSheet.Range("A$row")
    $aCell.Value = $d[0]
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $newSheet.Range("B$row").NumberFormat = "@"
    $newSheet.Range("B$row").Value = $d[1]
    $newSheet.Range("C$row").NumberFormat = "@"
    $newSheet.Range("C$row").Value = $d[2]
    $newSheet.Range("D$row").NumberFormat = "@"
    $newSheet.Range("D$row").Value = $d[3]
    $newSheet.Range("E$row").NumberFormat = "@"
    $newSheet.Range("E$row").Value = $d[4]
    $newSheet.Range("F$row").NumberFormat = "@"
    $newSheet.Range("F$row").Value = $d[5]
    $newSheet.Range("G$row").NumberFormat = "@"
    $newSheet.Range("G$row").Value = $d[6]
    $newSheet.Range("H$row").Value = $d[7]
    $row++
}
